# Edit LOQ4085.xlsx per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update Objetivos (row 10) B/C with new Portuguese objectives text
$ws.Range("B10").Value = 'Aplicar os fundamentos teóricos das operações unitárias envolvendo sistemas fluidos e particulados, baseados nos princípios dos fenômenos de transporte I.'
$ws.Range("C10").Value = 'Aplicar os fundamentos teóricos das operações unitárias envolvendo sistemas fluidos e particulados, baseados nos princípios dos fenômenos de transporte I.'

# 2) Insert a new row at 13 (shifts old rows 13-23 down to 14-24)
$ws.Rows("13:13").Insert()

# Clear the leftover A13 cell produced by the insert (target row 13 has no A cell)
$ws.Range("A13").Clear()

# Copy formats from the (now shifted) row 14 B/C cells so the new B13/C13
# pick up the correct column styles (content font + red-content font)
$ws.Range("B14").Copy() | Out-Null
$ws.Range("B13").PasteSpecial(-4122) | Out-Null
$ws.Range("C14").Copy() | Out-Null
$ws.Range("C13").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Populate new row 13 (Docentes responsaveis value)
$ws.Range("B13").Value = "8151869 - Livia Chaguri e Carvalho"
$ws.Range("C13").Value = "8151869 - Livia Chaguri e Carvalho"

# 3) Row 14 (was 13): Programa resumido -> now holds the short syllabus text
$ws.Range("B14").Value = "1)Transporte de fluidos (Newtonianos e não Newtonianos)2)Agitação e mistura3)Caracterização e dinâmica de partículas4)Separação de partículas por ação gravitacional e centrífuga5)Interação sólido – fluido6)Filtração7)Sedimentação"
$ws.Range("C14").Value = "1)Transporte de fluidos (Newtonianos e não Newtonianos)2)Agitação e mistura3)Caracterização e dinâmica de partículas4)Separação de partículas por ação gravitacional e centrífuga5)Interação sólido – fluido6)Filtração7)Sedimentação"

# 4) Row 16 (was 15 / "Programa:"): replace placeholder date with the full syllabus text
$ws.Range("B16").Value = '1)Transporte de fluidos: Tipos de bombas e compressores. Medidores de vazão. Curvas características. Cavitação e altura de sucção disponível (NPSH). Dimensionamento do sistema de bombeamento.2)Agitação e mistura: Tipos de equipamentos e impelidores. Mistura de líquidos. Cálculos de potência de agitadores.3)Caracterização e dinâmica de partículas: Características físicas de partícula isolada. Tamanho de partículas. Peneiramento. Análise granulométrica. Velocidade terminal.4)Separação de partículas por ação gravitacional e centrífuga: Elutriação. Câmara de poeira. Ciclones e centrífugas.5)Interação sólido – fluido: Escoamento em meio poroso. Fluidização.6)Filtração: Tipos de equipamentos. Filtração a pressão e vazão constante. Tortas compressíveis e incompressíveis.7)Sedimentação: Tipos de equipamentos. Cálculo da área e altura de sedimentadores.'
$ws.Range("C16").Value = '1)Transporte de fluidos: Tipos de bombas e compressores. Medidores de vazão. Curvas características. Cavitação e altura de sucção disponível (NPSH). Dimensionamento do sistema de bombeamento.2)Agitação e mistura: Tipos de equipamentos e impelidores. Mistura de líquidos. Cálculos de potência de agitadores.3)Caracterização e dinâmica de partículas: Características físicas de partícula isolada. Tamanho de partículas. Peneiramento. Análise granulométrica. Velocidade terminal.4)Separação de partículas por ação gravitacional e centrífuga: Elutriação. Câmara de poeira. Ciclones e centrífugas.5)Interação sólido – fluido: Escoamento em meio poroso. Fluidização.6)Filtração: Tipos de equipamentos. Filtração a pressão e vazão constante. Tortas compressíveis e incompressíveis.7)Sedimentação: Tipos de equipamentos. Cálculo da área e altura de sedimentadores.'

# 5) Row 19 (was 18 / "Método:"): replace placeholder teacher text with exam method text
$ws.Range("B19").Value = "Aplicação de 2 provas (P1 e P2)."
$ws.Range("C19").Value = "Aplicação de 2 provas (P1 e P2)."

# 6) Row 20 (was 19 / "Critério:"): was exam-method text, now the grading-average text
$ws.Range("B20").Value = "A média do período (MP) será calculada por: MP = (P1+P2)/2. `nAlunos com média final igual ou superior a 5,0 estarão aprovados, desde que tenham freqüência mínima de 70% (regimental). `nAlunos com média inferior a 3,0 e/ou freqüência inferior a 70% estarão reprovados (regimental). `nAlunos com média superior ou igual a 3,0 e inferior a 5,0 e que tenham freqüência mínima de 70% serão submetidos ao período de recuperação (regimental)."
$ws.Range("C20").Value = "A média do período (MP) será calculada por: MP = (P1+P2)/2. `nAlunos com média final igual ou superior a 5,0 estarão aprovados, desde que tenham freqüência mínima de 70% (regimental). `nAlunos com média inferior a 3,0 e/ou freqüência inferior a 70% estarão reprovados (regimental). `nAlunos com média superior ou igual a 3,0 e inferior a 5,0 e que tenham freqüência mínima de 70% serão submetidos ao período de recuperação (regimental)."

# 7) Row 21 (was 20 / "Norma de recuperação:"): was grading-average text, now recovery-average text; height 120 -> 60
$ws.Range("B21").Value = "A média final após a recuperação para a disciplina será a média aritmética entre a média do período e a nota da recuperação"
$ws.Range("C21").Value = "A média final após a recuperação para a disciplina será a média aritmética entre a média do período e a nota da recuperação"
$ws.Rows("21:21").RowHeight = 60

# 8) Row 22 (was 21 / "Bibliografia:"): was recovery-average text, now bibliography text; height none -> 120
$ws.Range("B22").Value = '1)COULSON, J. M.; RICHARDSON; J.F. Chemical Engineering. v.2: Particle Technology e Separation Processes. 5ed. Amsterdan: Butterworth Heinemann, 1229p. 2005;2)COULSON & Richardson''s Chemical Engineering: chemical engineering design by R.K. Sinnott. 6ed. Amsterdam: Elsevier Butterworth Heinemann, 895p. 2004;3)COUPER, J. R.; PENNEY, W. R.; FAIR, J. R.; W.; Stanley. M. Chemical Process Equipment: Selection and Design. 2ed. Amsterdam: Elsevier, 814p. 2005;4)MORAES JUNIOR, D. Transporte de líquidos e gases. v.1. São Carlos: Ufscar, 1988;5)FOUST, A. S.; WENZEL, L. A.; CLUMP, C. W.; MAUS, L.; ANDERSEN, L. B. 2ed. Princípios das operações unitárias. Rio de Janeiro: Guanabara Dois/LTC, 670p. 2008;6)GEANKOPLIS, C. J. Transport Processes and Separation Process Principles. 4ed. New York: Prentice Hall, 1026p. 2010;7)MCCABE, W. L.; SMITH, J. C.; HARRIOT, P. Unit operations of chemical engineering. 7ed. Boston: McGraw-Hill, 1140 p. 2005;8)PERRY''s chemical engineers handbook. Editor in Chief Don W. Green; Late Editor Robert H. Perry New York: McGraw-Hill, 2008.'
$ws.Range("C22").Value = '1)COULSON, J. M.; RICHARDSON; J.F. Chemical Engineering. v.2: Particle Technology e Separation Processes. 5ed. Amsterdan: Butterworth Heinemann, 1229p. 2005;2)COULSON & Richardson''s Chemical Engineering: chemical engineering design by R.K. Sinnott. 6ed. Amsterdam: Elsevier Butterworth Heinemann, 895p. 2004;3)COUPER, J. R.; PENNEY, W. R.; FAIR, J. R.; W.; Stanley. M. Chemical Process Equipment: Selection and Design. 2ed. Amsterdam: Elsevier, 814p. 2005;4)MORAES JUNIOR, D. Transporte de líquidos e gases. v.1. São Carlos: Ufscar, 1988;5)FOUST, A. S.; WENZEL, L. A.; CLUMP, C. W.; MAUS, L.; ANDERSEN, L. B. 2ed. Princípios das operações unitárias. Rio de Janeiro: Guanabara Dois/LTC, 670p. 2008;6)GEANKOPLIS, C. J. Transport Processes and Separation Process Principles. 4ed. New York: Prentice Hall, 1026p. 2010;7)MCCABE, W. L.; SMITH, J. C.; HARRIOT, P. Unit operations of chemical engineering. 7ed. Boston: McGraw-Hill, 1140 p. 2005;8)PERRY''s chemical engineers handbook. Editor in Chief Don W. Green; Late Editor Robert H. Perry New York: McGraw-Hill, 2008.'
$ws.Rows("22:22").RowHeight = 120

Write-Host "done"
